$wb = $excel.ActiveWorkbook

# Rename the worksheets (new timestamp-suffixed names)
$wb.Worksheets.Item(1).Name = "GNG_TO-16504777765843773"
$wb.Worksheets.Item(2).Name = "NB_TO-16504777807143772"
$wb.Worksheets.Item(3).Name = "RS_TO-1650477780715377"
$wb.Worksheets.Item(4).Name = "TOL_TO-165047778076239"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504777808254104"

# Sheet 1 (GNG_TO...) - update stimulus file names in column B
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650477776556376.csv"
$ws1.Range("B3").Value = "GNG_stims-16504777765674121.csv"
$ws1.Range("B4").Value = "go_stims-16504777765693796.csv"
$ws1.Range("B5").Value = "GNG_stims-16504777765834143.csv"

# Sheet 2 (NB_TO...) - update stimulus file names in column B
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16504777800294063.csv"
$ws2.Range("B3").Value = "ZB-match_4-1650477777278379.csv"
$ws2.Range("B4").Value = "TB-165047777949141.csv"
$ws2.Range("B5").Value = "OB-16504777775373774.csv"
$ws2.Range("B6").Value = "TB-16504777806953762.csv"
$ws2.Range("B7").Value = "ZB-match_9-16504777767244074.csv"
$ws2.Range("B8").Value = "OB-16504777793524096.csv"
$ws2.Range("B9").Value = "ZB-match_5-165047777716241.csv"
$ws2.Range("B10").Value = "OB-16504777785573783.csv"

# Sheet 3 (RS_TO...) - swap eyes closed / eyes open
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4 (TOL_TO...) - update stimulus file names in column B
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504777807293797.csv"
$ws4.Range("B3").Value = "ZM_stims-16504777807173786.csv"
$ws4.Range("B4").Value = "MM_stims-16504777807453783.csv"
$ws4.Range("B5").Value = "ZM_stims-16504777807293797.csv"
$ws4.Range("B6").Value = "MM_stims-16504777807613792.csv"
$ws4.Range("B7").Value = "ZM_stims-16504777807453783.csv"

# Sheet 5 (vSAT_TO...) - update stimulus file names in column B
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16504777807773817.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504777807933793.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504777808094099.csv"
$ws5.Range("B5").Value = "SAT_stims-16504777807653792.csv"
